$d = $word.ActiveDocument

# --- Change 1: remove the hidden "_GoBack" bookmark that used to sit on the
#     "For competitions, location may be unspecified" paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Change 2: drop the stale <w:lastRenderedPageBreak/> marker that sits
#     in front of the "Database Design" heading run. Re-assigning the
#     paragraph's own text rebuilds its run(s) without the marker while
#     preserving the run formatting (bold). ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq ("Database Design" + [char]13)) {
        $p.Range.Text = "Database Design"
    }
}

# --- Change 3: add the new walkthrough paragraphs describing how to use
#     inject_query, right after the "Database in Action" heading and before
#     the pre-existing (empty) paragraph that already followed it. The
#     "_GoBack" bookmark removed in Change 1 is re-created here, on its own
#     new trailing empty paragraph (this mirrors where Word's cursor last
#     was after the edit). ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq ("Database in Action" + [char]13)) {
        $target = $p
    }
}
$targetIndex = $target.Index
$insertPoint = $d.Range($target.Range.End, $target.Range.End)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="1140"/></w:tabs></w:pPr><w:r><w:t>Upon opening the bash shell,</w:t></w:r><w:r><w:t xml:space="preserve"> ensuring the project build is in the home folder</w:t></w:r><w:r><w:t xml:space="preserve">, and giving scripts in the bin folder execution </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>priveledges</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, we can initialize the database using</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="1140"/></w:tabs></w:pPr><w:r><w:tab/><w:t>“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>inject_query</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>yourusername</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>yourpassword</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>yourdatabase</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> db_init.txt”</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="1140"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>inject_query</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> script will log a user into the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mysql</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> database and run the specified query. Now, we should have tables in our database that hold all of our records.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="1140"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">We then navigate to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>db_queries</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> folder and, using the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>inject_query</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  script, run each query in order.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="1140"/></w:tabs></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>'

$insertPoint.InsertXML($xml) | Out-Null

# InsertXML above carries one extra, empty trailing <w:p/> so that Word
# splits our new content into its own paragraphs instead of merging the
# last of them into the pre-existing empty paragraph that used to directly
# follow "Database in Action". Remove that now-redundant placeholder
# paragraph (it lands right after our 5 new paragraphs).
$placeholder = $d.Paragraphs($targetIndex + 6)
$placeholder.Range.Delete() | Out-Null

Write-Output "done"
